$d = $word.ActiveDocument
$t3 = $d.Tables.Item(3)
$t3.Delete()
Write-Output ("Tables after delete: " + $d.Tables.Count)
Write-Output ("Paragraphs after delete: " + $d.Paragraphs.Count)
